$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.306.90'
$ws.Range('E2').Value = '  -4.59%  '

$ws.Range('D3').Value = '3.303.91'
$ws.Range('E3').Value = '  -0.52%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.78'
$ws.Range('E5').Value = '  -3.72%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.40'
$ws.Range('E6').Value = '  -5.71%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  -1.48%  '

$ws.Range('D9').Value = '3.304.20'
$ws.Range('E9').Value = '  -0.41%  '

$ws.Range('E10').Value = '  -3.94%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.62'
$ws.Range('E11').Value = '  -1.06%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.402'
$ws.Range('E12').Value = '  -4.47%  '

$ws.Range('D13').Value = '3.875.21'
$ws.Range('E13').Value = '  -0.57%  '

$ws.Range('E14').Value = '  -0.63%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.96'
$ws.Range('E15').Value = '  -4.43%  '

$ws.Range('D16').Value = '66.399.95'
$ws.Range('E16').Value = '  -4.41%  '

$ws.Range('E17').Value = '  -3.74%  '

$ws.Range('D18').Value = '3.326.32'
$ws.Range('E18').Value = '  -0.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '435.39'
$ws.Range('E19').Value = '  +3.28%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.61'
$ws.Range('E20').Value = '  -0.68%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.66'
$ws.Range('E21').Value = '  -2.93%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.63'
$ws.Range('E22').Value = '  -1.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.34'
$ws.Range('E23').Value = '  +1.66%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.14%  '

$ws.Range('D25').Value = '3.447.76'
$ws.Range('E25').Value = '  -0.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.511'
$ws.Range('E26').Value = '  -1.32%  '

$ws.Range('E27').Value = '  -3.13%  '

$ws.Range('E28').Value = '  -0.61%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.97'
$ws.Range('E29').Value = '  -7.34%  '

$ws.Range('E30').Value = '  -2.26%  '

$ws.Range('E31').Value = '  -1.72%  '

$ws.Range('E32').Value = '  -1.87%  '

$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.28'
$ws.Range('E34').Value = '  -6.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.74'
$ws.Range('E35').Value = '  -4.30%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.22'
$ws.Range('E36').Value = '  -5.13%  '

$ws.Range('E37').Value = '  -0.22%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.29'
$ws.Range('E38').Value = '  -2.85%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.20'
$ws.Range('E39').Value = '  +0.54%  '

$ws.Range('E40').Value = '  -5.55%  '

$ws.Range('D41').Value = '2.778.84'
$ws.Range('E41').Value = '  +2.43%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.783'
$ws.Range('E42').Value = '  -2.29%  '

$ws.Range('E43').Value = '  -3.33%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.16'
$ws.Range('E44').Value = '  -4.54%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0673'
$ws.Range('E45').Value = '  -2.64%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.14'
$ws.Range('E46').Value = '  -1.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.10'
$ws.Range('E47').Value = '  -5.33%  '

$ws.Range('E48').Value = '  -7.34%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '318.54'
$ws.Range('E49').Value = '  -7.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0271'
$ws.Range('E50').Value = '  -3.71%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.976'
$ws.Range('E51').Value = '  -3.27%  '
